$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8.428299999999995
$ws.Range("A3").Value = -22.07480000000001
$ws.Range("A14").Value = -21.8609
$ws.Range("A21").Value = -20.28309999999999
$ws.Range("A23").Value = -20.25849999999997
$ws.Range("A25").Value = -21.7339
$ws.Range("B25").Value = 5.493400000000001
$ws.Range("A26").Value = -21.02889999999997
$ws.Range("B27").Value = 6.140999999999999
$ws.Range("A29").Value = -20.63919999999998
$ws.Range("B31").Value = 5.557099999999999
$ws.Range("B39").Value = 9.610600000000003
$ws.Range("B48").Value = 5.187100000000003
$ws.Range("B51").Value = 5.542799999999996
$ws.Range("B52").Value = 5.273499999999998
$ws.Range("A53").Value = -21.9543
$ws.Range("B55").Value = 5.689799999999996
$ws.Range("B56").Value = 4.889899999999998
$ws.Range("A57").Value = -22.33980000000001
$ws.Range("B57").Value = 4.759999999999994
$ws.Range("A59").Value = -22.1306
$ws.Range("A69").Value = -21.62219999999999
$ws.Range("B73").Value = 8.440699999999994
$ws.Range("A79").Value = -20.2572
$ws.Range("A83").Value = -21.9952
$ws.Range("B89").Value = 4.641599999999995
$ws.Range("B90").Value = 5.350600000000002
$ws.Range("A91").Value = -21.46570000000002
$ws.Range("B92").Value = 4.765699999999999
$ws.Range("A93").Value = -20.71349999999998
